# Weekly update: insert a new (most recent) price record at the top of the
# data table (row 9) for "Vega Monumental Concepción - Poroto granado".
# Inserting the row shifts all existing data rows (old 9-53) down by one
# (new 10-54), which matches the target workbook exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 9; existing rows 9-53 shift down to 10-54.
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the latest observation.
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = "Vega Monumental Concepción"
$ws.Range("C9").Value = "Bíobío"
$ws.Range("D9").Value = 44950
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 100112030
$ws.Range("G9").Value = "Poroto granado"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 40000
$ws.Range("L9").Value = 42000
$ws.Range("M9").Value = 41000
$ws.Range("N9").Value = "`$/saco 25 kilos"
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value = 1640
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = "Hortaliza"
